{"js": "// Remove the \"Xian Gao\" paragraph from the CORE TEAM list in the\n// MEETING PARTICIPANTS table, leaving \"Molly Meadows\" as the last entry.\nconst hits = context.document.body.search(\"Xian Gao\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error('Could not find \"Xian Gao\" in the document body.');\n}\n\nconst hit = hits.items[0];\nconst paras = hit.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst para = paras.items[0];\npara.delete();\nawait context.sync();\n", "ps1": "# Remove the \"Xian Gao\" paragraph from the CORE TEAM list in the\n# MEETING PARTICIPANTS table, leaving \"Molly Meadows\" as the last entry.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if ($text -ne $null) {\n        $trimmed = $text.TrimEnd([char]13, [char]7)\n        if ($trimmed -eq \"Xian Gao\") {\n            $target = $p\n            break\n        }\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n"}
